$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: generation counter 16 -> 18 ---
$ws.Range("A1").Value = 18

# --- Extend the board two rows/cols (rows 4-21 get Q:R populated, board grows to R23) ---
$ws.Range("Q4:R21").Value = 0

# --- Row 20 is a brand-new all-zero board row ---
$ws.Range("A20:R20").Value = 0

# --- Row 21 used to hold only the "alive" flag (A21); now it's a plain zero board row ---
$ws.Range("A21:R21").Value = 0
$ws.Range("A21:R21").Style = "Normalny"

# --- The "alive" flag moves down to row 23 (row 22 stays blank, like row 3) ---
$ws.Range("A23").Value = "alive"
$ws.Range("A23").Interior.Color = 526344

# --- The live-cell ("glider") pattern on rows 9-12 shifts one column to the right ---
$clearCells = @("D9", "L9", "C10", "F10", "K10", "N10", "C11", "F11", "K11", "N11", "C12", "F12", "K12", "N12", "D13", "L13")
foreach ($c in $clearCells) {
    $ws.Range($c).Style = "Normalny"
}

$setCells = @("F9", "N9", "D10", "G10", "L10", "O10", "D11", "G11", "L11", "O11", "D12", "G12", "L12", "O12", "F13", "N13")
foreach ($c in $setCells) {
    $ws.Range($c).Interior.Color = 526344
}

# --- Selection moved to D2 ---
$ws.Range("D2").Select() | Out-Null
